$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Method row) updates
$ws.Range("C2").Value = 0.3182284544954106
$ws.Range("D2").Value = 0.3182284544954106
$ws.Range("E2").Value = 1.2472852808878423
$ws.Range("F2").Value = 0.008470684389926307
$ws.Range("G2").Value = 0.2762

# Row 3 (Residuals row) updates
$ws.Range("C3").Value = 37.24998207567866
$ws.Range("D3").Value = 0.2551368635320456
$ws.Range("F3").Value = 0.9915293156100736

# Row 4 (Total row) updates
$ws.Range("C4").Value = 37.56821053017407
